# Add 1 hour to the group time recorded for 17-02-2016 (cell D10, "week 2" column
# for the "gemeenschappelijk"/group row). The cell already holds a time-formatted
# value ([h]:mm), so we simply add 1/24 (one hour, expressed as a fraction of a day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D10")
$cell.Value2 = $cell.Value2 + (1 / 24)

# Move the active selection to D10, matching the saved view state in the file.
$ws.Range("D10").Select()
